$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text (Coin name) column B updates ---
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("B39").Value = 'MXToken'
$ws.Range("B40").Value = 'PaxDollar'
$ws.Range("B41").Value = 'VeChain'
$ws.Range("B42").Value = 'Quant'
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("B45").Value = 'Aptos'
$ws.Range("B46").Value = 'Algorand'
$ws.Range("B47").Value = 'Cronos'
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("B49").Value = 'Elrond'
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("B51").Value = 'Decentraland'

# --- Link column C updates ---
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("C40").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("C42").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("C45").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("C51").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'

# --- Price column D updates (force text to preserve exact formatting) ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.500.60'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.731.71'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9998'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.63'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4876'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.733.39'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.70'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.616'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6092'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '77.33'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.492.13'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007235'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.51'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.956.32'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.525'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.247'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '139.55'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.44'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.410'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '108.15'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.975'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08047'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.685'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04565'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.615'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.010'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6380'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9000'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.039'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.402'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.002'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.01513'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '101.34'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.439'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3889'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '6.966'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1184'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05395'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.840'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '30.56'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.249'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3407'

# --- Volume(1h) column E updates (force text to preserve exact formatting) ---
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.55%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.21%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +1.44%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.89%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.60%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.45%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.31%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.62%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.79%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.66%  '
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.02%  '
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.04%  '
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.29%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.04%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +4.82%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -2.05%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.52%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.71%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.44%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -2.11%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +2.72%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.55%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.66%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.88%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.64%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +2.23%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.47%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.02%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.16%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.90%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.47%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -3.74%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +3.15%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.38%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.31%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.02%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -10.46%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -5.44%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.50%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +2.84%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.92%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.18%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.42%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.63%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.80%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.98%  '
